{"js": "// 1. Remove the stray \"_GoBack\" bookmark that sits at the very start of the\n//    document (before the title text \"Wat moet er nog gebeuren voor animate\").\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. The \"Hardware en onderdelen net zo iets als die van planty\" bullet\n//    (currently the 4th paragraph) now needs to be struck through, like the\n//    \"Fix home button\" bullet already is.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst hardwareParagraph = paragraphs.items[3];\nhardwareParagraph.font.strikeThrough = true;\n\n// 3. Word's \"last edit position\" bookmark (\"_GoBack\") now belongs in the\n//    final paragraph, right after \"...der aan Casper als\" and before\n//    \" hij opmerkingen heeft over wat misschien makkelijk beter kan \".\nconst searchResults = body.search(\"der aan Casper als\", { matchCase: false });\nsearchResults.load(\"text\");\nawait context.sync();\n\nconst matchRange = searchResults.items[0];\nconst insertionPoint = matchRange.getRange(\"After\");\ninsertionPoint.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the stray \"_GoBack\" bookmark that sits at the very start of the\n#    document (before the title text \"Wat moet er nog gebeuren voor animate\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. The \"Hardware en onderdelen net zo iets als die van planty\" bullet\n#    (the 4th paragraph) now needs to be struck through, like the\n#    \"Fix home button\" bullet already is.\n$hardwareParagraph = $d.Paragraphs.Item(4)\n$hardwareParagraph.Range.Font.StrikeThrough = 1\n\n# 3. Word's \"last edit position\" bookmark (\"_GoBack\") now belongs in the\n#    final paragraph, right after \"...der aan Casper als\" and before\n#    \" hij opmerkingen heeft over wat misschien makkelijk beter kan \".\n$findRange = $d.Content\n$null = $findRange.Find.Execute(\"der aan Casper als\")\n$insertionPoint = $d.Range($findRange.End, $findRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n"}
